# Adds a new "2022-Q4" sheet (positioned right after "总计") with its
# fund-holding breakdown, and updates the "总计" (summary) sheet so its
# first data row reflects the new 2022-Q4 totals while the previously
# existing quarters shift down by one row. The other quarter sheets
# (2022-Q3 / 2022-Q2 / 2022-Q1 / 2021-Q4) are left untouched and simply
# slide one tab to the right because the new sheet is inserted before them.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Helper: write a value as literal TEXT (matches source cells such as
# "34.68" / "4.00" which are stored as strings, not numbers).
# ---------------------------------------------------------------------
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Header row (copy style/format from the "总计" sheet header so the
# bold + bordered + centered look matches the other quarter sheets).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2 # column B
foreach ($h in $headers) {
    $cell = $q4.Cells.Item(1, $col)
    $cell.Value = $h
    $total.Range("B1").Copy()
    $cell.PasteSpecial(-4122)
    $col = $col + 1
}

# Data rows for the new quarter.
$q4rows = @(
    @{A=0; B="011479"; C="广发诚享混合A";             D="34.68";  E="93.23"; F="5.21"; G="1.8068"; H=9},
    @{A=1; B="011130"; C="广发兴诚混合C";             D="24.57";  E="94.38"; F="4.56"; G="1.1204"; H=10},
    @{A=2; B="011121"; C="广发兴诚混合A";             D="21.01";  E="94.38"; F="4.56"; G="0.9581"; H=10},
    @{A=3; B="011480"; C="广发诚享混合C";             D="4.00";   E="93.23"; F="5.21"; G="0.2084"; H=9},
    @{A=4; B="004044"; C="金鹰转型动力灵活配置混合"; D="0.65";   E="90.51"; F="4.09"; G="0.0266"; H=10}
)

$r = 2
foreach ($row in $q4rows) {
    $cellA = $q4.Cells.Item($r, 1)
    $cellA.Value = $row.A
    $total.Range("A2").Copy()
    $cellA.PasteSpecial(-4122)
    $cellA.Value = $row.A

    $q4.Cells.Item($r, 2).Value = $row.B
    $q4.Cells.Item($r, 3).Value = $row.C
    Set-TextValue $q4.Cells.Item($r, 4) $row.D
    Set-TextValue $q4.Cells.Item($r, 5) $row.E
    Set-TextValue $q4.Cells.Item($r, 6) $row.F
    Set-TextValue $q4.Cells.Item($r, 7) $row.G
    $q4.Cells.Item($r, 8).Value = $row.H

    $r = $r + 1
}

$q4.Range("A1").Select()

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: row 2 becomes the new 2022-Q4 entry, and
#    the rows that used to hold 2022-Q3 / 2022-Q2 / 2022-Q1 / 2021-Q4
#    shift down by one (row 6 is brand new, so its column-A style is
#    copied from row 5 first).
# ---------------------------------------------------------------------
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)

$totalRows = @(
    @{Row=2; A=0; B="2022-Q4"; C=5;  D=4.12},
    @{Row=3; A=1; B="2022-Q3"; C=23; D=2.75},
    @{Row=4; A=2; B="2022-Q2"; C=8;  D=7.63},
    @{Row=5; A=3; B="2022-Q1"; C=8;  D=9.46},
    @{Row=6; A=4; B="2021-Q4"; C=4;  D=10.05}
)

foreach ($row in $totalRows) {
    $total.Cells.Item($row.Row, 1).Value = $row.A
    $total.Cells.Item($row.Row, 2).Value = $row.B
    $total.Cells.Item($row.Row, 3).Value = $row.C
    $total.Cells.Item($row.Row, 4).Value = $row.D
}

$total.Activate()
$total.Range("A1").Select()
